# Adds the missing MIAPPE ontology reference (MIAPPE:0166) to the
# "Factor [other perturbation]" building block of the MIAPPE experimental
# factors template, and bumps the template version 1.0.2 -> 1.0.3.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")

# --- Rename the three trailing "Factor [other perturbation]" columns of the
#     annotation table so they carry the MIAPPE:0166 term accession, fixing
#     the casing of "other" -> "Other" to match the other factor headers.
$ws1.Range("BV1").Value = "Factor [Other perturbation]"
$ws1.Range("BW1").Value = "Term Source REF (MIAPPE:0166)"
$ws1.Range("BX1").Value = "Term Accession Number (MIAPPE:0166)"

# --- The longer header text needs slightly wider (autofit) columns, matching
#     the width used by all the other Term Source REF / Term Accession
#     Number column pairs in the sheet.
$ws1.Range("BV1").EntireColumn.ColumnWidth = 25.43
$ws1.Range("BW1").EntireColumn.ColumnWidth = 29.83
$ws1.Range("BX1").EntireColumn.ColumnWidth = 36.33

# --- Bump the template version.
$ws2.Range("B3").Value = "1.0.3"

# --- Restore view state: Tabelle1 keeps its scroll/selection from before it
#     was left, SwateTemplateMetadata (tab 2) stays the active sheet/tab.
$ws1.Activate()
$ws1.Range("BV6").Select()

$ws2.Activate()
$ws2.Range("B4").Select()
